# Insert a new row at row 540, shifting existing rows 540:578 down to 541:579,
# then populate the newly inserted row 540 with the new data record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before current row 540 (Excel semantics: shifts down).
$ws.Rows.Item(540).Insert()

# Populate the new row 540 with the data for the new weekly record.
$ws.Range("A540").Value = 10
$ws.Range("B540").Value = "Vega Modelo de Temuco"
$ws.Range("C540").Value = "La Araucanía"
$ws.Range("D540").Value = 44826
$ws.Range("E540").Value = 9
$ws.Range("F540").Value = 100112043
$ws.Range("G540").Value = "Pepino ensalada"
$ws.Range("H540").Value = "Sin especificar"
$ws.Range("I540").Value = "Primera"
$ws.Range("J540").Value = 500
$ws.Range("K540").Value = 20000
$ws.Range("L540").Value = 24000
$ws.Range("M540").Value = 22400
$ws.Range("N540").Value = '$/caja 60 unidades'
$ws.Range("O540").Value = "Región de Arica y Parinacota"
$ws.Range("P540").Value = 373
$ws.Range("Q540").Value = 60
$ws.Range("R540").Value = "Hortaliza"
